$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A36").Value = "test"
